$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number-format/border/bold style used by the existing column-A
# "index" cells (A3:A16) onto the three new index cells before filling them
# in, so the new rows match the look of the existing ones.
$ws.Range("A3").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

# New rows 17-19: averaged-intensity results for the spiral sampling
# schemes that were just run.
$rowData = @(
    @{ Row = 17; Idx = 15; Label = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; Idx = 16; Label = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; Idx = 17; Label = "HexGrid-60degTilt5degRes" }
)

foreach ($r in $rowData) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Idx
    $ws.Range("B$rowNum").Value = $r.Label
    foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M")) {
        $ws.Range("$col$rowNum").Value = 1
    }
}
